# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-numeric-looking price strings to be stored as text (matches source formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values
# Row 2
$ws.Range("D2").Value = "40.901.46"
$ws.Range("E2").Value = "  -1.78%  "
# Row 3
$ws.Range("D3").Value = "2.421.48"
$ws.Range("E3").Value = "  -2.12%  "
# Row 4
$ws.Range("E4").Value = "  +0.11%  "
# Row 5
$ws.Range("D5").Value = "316.65"
$ws.Range("E5").Value = "  -0.47%  "
# Row 6
$ws.Range("D6").Value = "89.16"
$ws.Range("E6").Value = "  -3.53%  "
# Row 7
$ws.Range("D7").Value = "0.538"
$ws.Range("E7").Value = "  -2.84%  "
# Row 8
$ws.Range("E8").Value = "  +0.05%  "
# Row 9
$ws.Range("D9").Value = "0.496"
$ws.Range("E9").Value = "  -3.92%  "
# Row 10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.0832"
$ws.Range("E10").Value = "  -4.09%  "
# Row 11
$ws.Range("B11").Value = "Avalanche"
$ws.Range("C11").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D11").Value = "31.81"
$ws.Range("E11").Value = "  -4.04%  "
# Row 12
$ws.Range("E12").Value = "  -2.34%  "
# Row 13
$ws.Range("D13").Value = "2.790.05"
$ws.Range("E13").Value = "  -2.26%  "
# Row 14
$ws.Range("D14").Value = "6.70"
$ws.Range("E14").Value = "  -2.90%  "
# Row 15
$ws.Range("D15").Value = "15.77"
$ws.Range("E15").Value = "  +0.80%  "
# Row 16
$ws.Range("D16").Value = "2.399.47"
$ws.Range("E16").Value = "  -2.90%  "
# Row 17
$ws.Range("D17").Value = "0.769"
$ws.Range("E17").Value = "  -2.84%  "
# Row 18
$ws.Range("D18").Value = "40.857.36"
$ws.Range("E18").Value = "  -1.78%  "
# Row 19
$ws.Range("D19").Value = "0.0₃0924"
$ws.Range("E19").Value = "  -3.21%  "
# Row 20
$ws.Range("D20").Value = "6.25"
$ws.Range("E20").Value = "  -3.84%  "
# Row 21
$ws.Range("D21").Value = "71.27"
$ws.Range("E21").Value = "  -0.18%  "
# Row 22
$ws.Range("D22").Value = "11.01"
$ws.Range("E22").Value = "  -3.11%  "
# Row 23
$ws.Range("D23").Value = "235.00"
$ws.Range("E23").Value = "  -2.48%  "
# Row 24
$ws.Range("E24").Value = "  -2.34%  "
# Row 25
$ws.Range("E25").Value = "  +0.14%  "
# Row 26
$ws.Range("D26").Value = "1.89"
$ws.Range("E26").Value = "  -2.27%  "
# Row 27
$ws.Range("D27").Value = "24.12"
$ws.Range("E27").Value = "  -2.65%  "
# Row 28
$ws.Range("D28").Value = "2.22"
$ws.Range("E28").Value = "  -3.05%  "
# Row 29
$ws.Range("D29").Value = "9.55"
$ws.Range("E29").Value = "  -3.45%  "
# Row 30
$ws.Range("D30").Value = "34.82"
$ws.Range("E30").Value = "  -4.09%  "
# Row 31
$ws.Range("D31").Value = "155.35"
$ws.Range("E31").Value = "  -2.83%  "
# Row 32
$ws.Range("E32").Value = "  -0.01%  "
# Row 33
$ws.Range("D33").Value = "5.26"
$ws.Range("E33").Value = "  -4.96%  "
# Row 34
$ws.Range("E34").Value = "  -2.80%  "
# Row 35
$ws.Range("D35").Value = "0.0743"
$ws.Range("E35").Value = "  -3.93%  "
# Row 36
$ws.Range("D36").Value = "2.98"
$ws.Range("E36").Value = "  +2.59%  "
# Row 37
$ws.Range("D37").Value = "16.70"
$ws.Range("E37").Value = "  -3.97%  "
# Row 38
$ws.Range("D38").Value = "0.115"
$ws.Range("E38").Value = "  -0.91%  "
# Row 39
$ws.Range("D39").Value = "1.78"
$ws.Range("E39").Value = "  -2.85%  "
# Row 40
$ws.Range("D40").Value = "0.100"
$ws.Range("E40").Value = "  -2.73%  "
# Row 41
$ws.Range("D41").Value = "3.89"
$ws.Range("E41").Value = "  -2.10%  "
# Row 42
$ws.Range("D42").Value = "1.990.86"
$ws.Range("E42").Value = "  +0.04%  "
# Row 43
$ws.Range("E43").Value = "  -8.36%  "
# Row 44
$ws.Range("D44").Value = "19.11"
$ws.Range("E44").Value = "  -0.36%  "
# Row 45
$ws.Range("D45").Value = "0.0275"
$ws.Range("E45").Value = "  -4.03%  "
# Row 46
$ws.Range("D46").Value = "2.88"
$ws.Range("E46").Value = "  -4.01%  "
# Row 47
$ws.Range("D47").Value = "9.51"
$ws.Range("E47").Value = "  +3.44%  "
# Row 48
$ws.Range("D48").Value = "2.650.77"
$ws.Range("E48").Value = "  -2.25%  "
# Row 49
$ws.Range("D49").Value = "94.87"
$ws.Range("E49").Value = "  -2.98%  "
# Row 50
$ws.Range("D50").Value = "73.41"
$ws.Range("E50").Value = "  -0.44%  "
# Row 51
$ws.Range("D51").Value = "51.93"
$ws.Range("E51").Value = "  -1.15%  "
